$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.345.31'
$ws.Range("E2").Value = '  -0.56%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.307.08'
$ws.Range("E3").Value = '  -0.75%  '

$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '518.25'
$ws.Range("E5").Value = '  +0.28%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.45'
$ws.Range("E6").Value = '  -3.15%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.530'
$ws.Range("E8").Value = '  -1.35%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.317.45'
$ws.Range("E9").Value = '  -1.11%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0997'
$ws.Range("E10").Value = '  -2.66%  '

$ws.Range("E11").Value = '  +0.07%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.25'
$ws.Range("E12").Value = '  -1.76%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.337'
$ws.Range("E13").Value = '  -1.91%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.718.93'
$ws.Range("E14").Value = '  -0.79%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.32'
$ws.Range("E15").Value = '  -2.69%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '56.302.81'
$ws.Range("E16").Value = '  -0.73%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000132'
$ws.Range("E17").Value = '  -2.19%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.321.03'
$ws.Range("E18").Value = '  -0.44%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '331.40'
$ws.Range("E19").Value = '  +1.47%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.33'
$ws.Range("E20").Value = '  -2.02%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.13'
$ws.Range("E21").Value = '  -2.32%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.70'
$ws.Range("E22").Value = '  +1.57%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.97'
$ws.Range("E24").Value = '  +0.30%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.62'
$ws.Range("E25").Value = '  +7.81%  '

$ws.Range("E26").Value = '  -0.79%  '

$ws.Range("E27").Value = '  +0.04%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.32'
$ws.Range("E28").Value = '  +2.66%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '168.60'
$ws.Range("E29").Value = '  -0.38%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.69'
$ws.Range("E30").Value = '  -0.33%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0₃0714'
$ws.Range("E31").Value = '  -4.10%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.08'
$ws.Range("E32").Value = '  -2.00%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.21'
$ws.Range("E33").Value = '  -1.45%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.992'
$ws.Range("E35").Value = '  -0.06%  '

$ws.Range("E36").Value = '  -2.02%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.90'
$ws.Range("E37").Value = '  -2.65%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.881'
$ws.Range("E38").Value = '  -4.08%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.57'
$ws.Range("E39").Value = '  +0.36%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '38.55'
$ws.Range("E40").Value = '  +0.55%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '148.30'
$ws.Range("E41").Value = '  +4.26%  '

$ws.Range("E42").Value = '  -2.51%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '284.09'
$ws.Range("E43").Value = '  +1.97%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.55'
$ws.Range("E44").Value = '  -1.38%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.05'
$ws.Range("E45").Value = '  -2.68%  '

$ws.Range("E46").Value = '  -1.32%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0496'
$ws.Range("E47").Value = '  -2.14%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '18.18'
$ws.Range("E49").Value = '  +1.57%  '

$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0213'
$ws.Range("E50").Value = '  -2.55%  '

$ws.Range("B51").Value = 'Polygon'
$ws.Range("C51").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.376'
$ws.Range("E51").Value = '  -1.07%  '
